$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the measured angle values for specimen 4 (row 5)
$ws.Range("B5").Value = 17.515
$ws.Range("C5").Value = 17.707
$ws.Range("D5").Value = 21.404

# Fill the averaging formula down into E5 (shared formula group already spans E3:E11)
$ws.Range("E5").Formula = "=AVERAGE(C5:D5)"

# Update the active selection to match the edited range
$ws.Range("E4:E5").Select()
